# Apply the marksheet corrections:
# - B11 (Marking/Right count): 3 -> 5
# - B12 (Total/Right count):   69 -> 115
# - E12 (Total/Max text):      "68/84" -> "115/140"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 115
$ws.Range("E12").Value = "115/140"
